$wb = $excel.ActiveWorkbook

# ---- Sheet index 16 ----
$ws = $wb.Worksheets.Item(16)
$ws.Range("A4:N4").Copy($ws.Range("A21:N21"))
$ws.Range("A21").Value = "EL PRIMO"
$ws.Range("B21").Value = "MELODIE"
$ws.Range("C21").Value = "SHADE"
$ws.Range("D21").Value = "KAZE"
$ws.Range("E21").Value = "HANK"
$ws.Range("F21").Value = "LOU"
$ws.Range("G21").Value = "Equipo 1"
$ws.Range("H21").Value = "PLP|BrriN"
$ws.Range("I21").Value = "MTM|snoiy"
$ws.Range("J21").Value = "PLP|Mine"
$ws.Range("K21").Value = "NHG|Xemp"
$ws.Range("L21").Value = "NHG|GN"
$ws.Range("M21").Value = "NHG|Bayarea"
$ws.Range("N21").Value = "20250724T012024.000Z"
$ws.Range("A20:N20").Copy($ws.Range("A22:N22"))
$ws.Range("A22").Value = "EL PRIMO"
$ws.Range("B22").Value = "MELODIE"
$ws.Range("C22").Value = "SHADE"
$ws.Range("D22").Value = "KAZE"
$ws.Range("E22").Value = "HANK"
$ws.Range("F22").Value = "LOU"
$ws.Range("G22").Value = "Equipo 2"
$ws.Range("H22").Value = "PLP|BrriN"
$ws.Range("I22").Value = "MTM|snoiy"
$ws.Range("J22").Value = "PLP|Mine"
$ws.Range("K22").Value = "NHG|Xemp"
$ws.Range("L22").Value = "NHG|GN"
$ws.Range("M22").Value = "NHG|Bayarea"
$ws.Range("N22").Value = "20250724T011836.000Z"
$ws.Range("A20:N20").Copy($ws.Range("A23:N23"))
$ws.Range("A23").Value = "BULL"
$ws.Range("B23").Value = "BERRY"
$ws.Range("C23").Value = "LOU"
$ws.Range("D23").Value = "WILLOW"
$ws.Range("E23").Value = "KAZE"
$ws.Range("F23").Value = "HANK"
$ws.Range("G23").Value = "Equipo 2"
$ws.Range("H23").Value = "PLP|BrriN"
$ws.Range("I23").Value = "MTM|snoiy"
$ws.Range("J23").Value = "PLP|Mine"
$ws.Range("K23").Value = "NHG|Bayarea"
$ws.Range("L23").Value = "NHG|GN"
$ws.Range("M23").Value = "NHG|Xemp"
$ws.Range("N23").Value = "20250724T011227.000Z"
$ws.Range("A20:N20").Copy($ws.Range("A24:N24"))
$ws.Range("A24").Value = "BULL"
$ws.Range("B24").Value = "BERRY"
$ws.Range("C24").Value = "LOU"
$ws.Range("D24").Value = "WILLOW"
$ws.Range("E24").Value = "KAZE"
$ws.Range("F24").Value = "HANK"
$ws.Range("G24").Value = "Equipo 2"
$ws.Range("H24").Value = "PLP|BrriN"
$ws.Range("I24").Value = "MTM|snoiy"
$ws.Range("J24").Value = "PLP|Mine"
$ws.Range("K24").Value = "NHG|Bayarea"
$ws.Range("L24").Value = "NHG|GN"
$ws.Range("M24").Value = "NHG|Xemp"
$ws.Range("N24").Value = "20250724T011025.000Z"
$ws.Range("A4:N4").Copy($ws.Range("A25:N25"))
$ws.Range("A25").Value = "HANK"
$ws.Range("B25").Value = "KAZE"
$ws.Range("C25").Value = "NITA"
$ws.Range("D25").Value = "LUMI"
$ws.Range("E25").Value = "OTIS"
$ws.Range("F25").Value = "MICO"
$ws.Range("G25").Value = "Equipo 1"
$ws.Range("H25").Value = "PLP|BrriN"
$ws.Range("I25").Value = "MTM|snoiy"
$ws.Range("J25").Value = "PLP|Mine"
$ws.Range("K25").Value = "NHG|GN"
$ws.Range("L25").Value = "NHG|Bayarea"
$ws.Range("M25").Value = "NHG|Xemp"
$ws.Range("N25").Value = "20250724T010505.000Z"
$ws.Range("A4:N4").Copy($ws.Range("A26:N26"))
$ws.Range("A26").Value = "HANK"
$ws.Range("B26").Value = "KAZE"
$ws.Range("C26").Value = "NITA"
$ws.Range("D26").Value = "LUMI"
$ws.Range("E26").Value = "OTIS"
$ws.Range("F26").Value = "MICO"
$ws.Range("G26").Value = "Equipo 1"
$ws.Range("H26").Value = "PLP|BrriN"
$ws.Range("I26").Value = "MTM|snoiy"
$ws.Range("J26").Value = "PLP|Mine"
$ws.Range("K26").Value = "NHG|GN"
$ws.Range("L26").Value = "NHG|Bayarea"
$ws.Range("M26").Value = "NHG|Xemp"
$ws.Range("N26").Value = "20250724T010329.000Z"
$ws.Range("A20:N20").Copy($ws.Range("A27:N27"))
$ws.Range("A27").Value = "HANK"
$ws.Range("B27").Value = "KAZE"
$ws.Range("C27").Value = "NITA"
$ws.Range("D27").Value = "LUMI"
$ws.Range("E27").Value = "OTIS"
$ws.Range("F27").Value = "MICO"
$ws.Range("G27").Value = "Equipo 2"
$ws.Range("H27").Value = "PLP|BrriN"
$ws.Range("I27").Value = "MTM|snoiy"
$ws.Range("J27").Value = "PLP|Mine"
$ws.Range("K27").Value = "NHG|GN"
$ws.Range("L27").Value = "NHG|Bayarea"
$ws.Range("M27").Value = "NHG|Xemp"
$ws.Range("N27").Value = "20250724T010039.000Z"

# ---- Sheet index 4 ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("A47:N47").Copy($ws.Range("A48:N48"))
$ws.Range("A48").Value = "AMBER"
$ws.Range("B48").Value = "BERRY"
$ws.Range("C48").Value = "LILY"
$ws.Range("D48").Value = "RICO"
$ws.Range("E48").Value = "DRACO"
$ws.Range("F48").Value = "KAZE"
$ws.Range("G48").Value = "Equipo 1"
$ws.Range("H48").Value = "TE|Rafikii"
$ws.Range("I48").Value = "TE|Ezlivi"
$ws.Range("J48").Value = "TE|Belal"
$ws.Range("K48").Value = "TRB|Zeus 解開"
$ws.Range("L48").Value = "TRB|R B M"
$ws.Range("M48").Value = "TRB|Lxffy"
$ws.Range("N48").Value = "20250724T011858.000Z"
$ws.Range("A47:N47").Copy($ws.Range("A49:N49"))
$ws.Range("A49").Value = "AMBER"
$ws.Range("B49").Value = "BERRY"
$ws.Range("C49").Value = "LILY"
$ws.Range("D49").Value = "RICO"
$ws.Range("E49").Value = "DRACO"
$ws.Range("F49").Value = "KAZE"
$ws.Range("G49").Value = "Equipo 1"
$ws.Range("H49").Value = "TE|Rafikii"
$ws.Range("I49").Value = "TE|Ezlivi"
$ws.Range("J49").Value = "TE|Belal"
$ws.Range("K49").Value = "TRB|Zeus 解開"
$ws.Range("L49").Value = "TRB|R B M"
$ws.Range("M49").Value = "TRB|Lxffy"
$ws.Range("N49").Value = "20250724T011722.000Z"
$ws.Range("A4:N4").Copy($ws.Range("A50:N50"))
$ws.Range("A50").Value = "AMBER"
$ws.Range("B50").Value = "BONNIE"
$ws.Range("C50").Value = "BERRY"
$ws.Range("D50").Value = "GRIFF"
$ws.Range("E50").Value = "SAM"
$ws.Range("F50").Value = "KAZE"
$ws.Range("G50").Value = "Equipo 2"
$ws.Range("H50").Value = "TE|Rafikii"
$ws.Range("I50").Value = "TE|Ezlivi"
$ws.Range("J50").Value = "TE|Belal"
$ws.Range("K50").Value = "TRB|Zeus 解開"
$ws.Range("L50").Value = "T1|Keria"
$ws.Range("M50").Value = "TRB|Lxffy"
$ws.Range("N50").Value = "20250724T011146.000Z"
$ws.Range("A4:N4").Copy($ws.Range("A51:N51"))
$ws.Range("A51").Value = "AMBER"
$ws.Range("B51").Value = "BONNIE"
$ws.Range("C51").Value = "BERRY"
$ws.Range("D51").Value = "GRIFF"
$ws.Range("E51").Value = "DARRYL"
$ws.Range("F51").Value = "KAZE"
$ws.Range("G51").Value = "Equipo 2"
$ws.Range("H51").Value = "TE|Rafikii"
$ws.Range("I51").Value = "TE|Ezlivi"
$ws.Range("J51").Value = "TE|Belal"
$ws.Range("K51").Value = "TRB|Zeus 解開"
$ws.Range("L51").Value = "T1|Keria"
$ws.Range("M51").Value = "TRB|Lxffy"
$ws.Range("N51").Value = "20250724T010850.000Z"
$ws.Range("A4:N4").Copy($ws.Range("A52:N52"))
$ws.Range("A52").Value = "AMBER"
$ws.Range("B52").Value = "BONNIE"
$ws.Range("C52").Value = "BERRY"
$ws.Range("D52").Value = "GRIFF"
$ws.Range("E52").Value = "DARRYL"
$ws.Range("F52").Value = "KAZE"
$ws.Range("G52").Value = "Equipo 2"
$ws.Range("H52").Value = "TE|Rafikii"
$ws.Range("I52").Value = "TE|Ezlivi"
$ws.Range("J52").Value = "TE|Belal"
$ws.Range("K52").Value = "TRB|Zeus 解開"
$ws.Range("L52").Value = "T1|Keria"
$ws.Range("M52").Value = "TRB|Lxffy"
$ws.Range("N52").Value = "20250724T010717.000Z"

# ---- Sheet index 5 ----
$ws = $wb.Worksheets.Item(5)
$ws.Range("A49:N49").Copy($ws.Range("A50:N50"))
$ws.Range("A50").Value = "JUJU"
$ws.Range("B50").Value = "MEEPLE"
$ws.Range("C50").Value = "HANK"
$ws.Range("D50").Value = "GENE"
$ws.Range("E50").Value = "MICO"
$ws.Range("F50").Value = "CORDELIUS"
$ws.Range("G50").Value = "Equipo 1"
$ws.Range("H50").Value = "TE|Rafikii"
$ws.Range("I50").Value = "TE|Ezlivi"
$ws.Range("J50").Value = "TE|Belal"
$ws.Range("K50").Value = "TRB|Zeus 解開"
$ws.Range("L50").Value = "TRB|R B M"
$ws.Range("M50").Value = "TRB|Lxffy"
$ws.Range("N50").Value = "20250724T010145.000Z"
$ws.Range("A4:N4").Copy($ws.Range("A51:N51"))
$ws.Range("A51").Value = "JUJU"
$ws.Range("B51").Value = "MEEPLE"
$ws.Range("C51").Value = "HANK"
$ws.Range("D51").Value = "GENE"
$ws.Range("E51").Value = "MICO"
$ws.Range("F51").Value = "CORDELIUS"
$ws.Range("G51").Value = "Equipo 2"
$ws.Range("H51").Value = "TE|Rafikii"
$ws.Range("I51").Value = "TE|Ezlivi"
$ws.Range("J51").Value = "TE|Belal"
$ws.Range("K51").Value = "TRB|Zeus 解開"
$ws.Range("L51").Value = "TRB|R B M"
$ws.Range("M51").Value = "TRB|Lxffy"
$ws.Range("N51").Value = "20250724T005945.000Z"
$ws.Range("A49:N49").Copy($ws.Range("A52:N52"))
$ws.Range("A52").Value = "JUJU"
$ws.Range("B52").Value = "MEEPLE"
$ws.Range("C52").Value = "HANK"
$ws.Range("D52").Value = "GENE"
$ws.Range("E52").Value = "MICO"
$ws.Range("F52").Value = "CORDELIUS"
$ws.Range("G52").Value = "Equipo 1"
$ws.Range("H52").Value = "TE|Rafikii"
$ws.Range("I52").Value = "TE|Ezlivi"
$ws.Range("J52").Value = "TE|Belal"
$ws.Range("K52").Value = "TRB|Zeus 解開"
$ws.Range("L52").Value = "TRB|R B M"
$ws.Range("M52").Value = "TRB|Lxffy"
$ws.Range("N52").Value = "20250724T005725.000Z"
$ws.Range("A49:N49").Copy($ws.Range("A53:N53"))
$ws.Range("A53").Value = "DOUG"
$ws.Range("B53").Value = "PENNY"
$ws.Range("C53").Value = "JANET"
$ws.Range("D53").Value = "GUS"
$ws.Range("E53").Value = "MR. P"
$ws.Range("F53").Value = "CORDELIUS"
$ws.Range("G53").Value = "Equipo 1"
$ws.Range("H53").Value = "TE|Rafikii"
$ws.Range("I53").Value = "TE|Ezlivi"
$ws.Range("J53").Value = "TE|Belal"
$ws.Range("K53").Value = "TRB|Zeus 解開"
$ws.Range("L53").Value = "TRB|R B M"
$ws.Range("M53").Value = "TRB|Lxffy"
$ws.Range("N53").Value = "20250724T005117.000Z"
